$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2..306 was updated from serial date
# 45186 (2023-09-17) to 45188 (2023-09-19).
$ws.Range("C2:C306").Value = 45188
